$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): replace "/Parent/Child" style labels with "Parent-Child" style ---
$ws.Range("B1").Value = "Government-Cadw"
$ws.Range("C1").Value = "Government-Local_Authority"
$ws.Range("D1").Value = "Government-National"
$ws.Range("E1").Value = "Government-Other"
$ws.Range("F1").Value = "Independent-English_Heritage"
$ws.Range("G1").Value = "Independent-Historic_Environment_Scotland"
$ws.Range("H1").Value = "Independent-National_Trust"
$ws.Range("I1").Value = "Independent-National_Trust_for_Scotland"
$ws.Range("J1").Value = "Independent-Not_for_profit"
$ws.Range("K1").Value = "Independent-Private"
$ws.Range("L1").Value = "Independent-Unknown"
$ws.Range("M1").Value = "University"
$ws.Range("N1").Value = "Unknown"

# --- Data rows 2-7: updated statistics values ---

# Row 2 (small)
$ws.Range("B2").Value = 0.071
$ws.Range("C2").Value = 6.432
$ws.Range("D2").Value = 0.142
$ws.Range("E2").Value = 0.047
$ws.Range("F2").Value = 0.284
$ws.Range("G2").Value = 0.166
$ws.Range("H2").Value = 0.591
$ws.Range("I2").Value = 0.166
$ws.Range("J2").Value = 27.903
$ws.Range("K2").Value = 13.951
$ws.Range("L2").Value = 4.611
$ws.Range("M2").Value = 1.088
$ws.Range("N2").Value = 1.939
$ws.Range("O2").Value = 57.391

# Row 3 (medium)
$ws.Range("C3").Value = 10.641
$ws.Range("D3").Value = 0.213
$ws.Range("E3").Value = 0.071
$ws.Range("F3").Value = 0.638
$ws.Range("G3").Value = 0.189
$ws.Range("H3").Value = 1.561
$ws.Range("I3").Value = 0.355
$ws.Range("J3").Value = 7.945
$ws.Range("K3").Value = 1.182
$ws.Range("L3").Value = 0.402
$ws.Range("M3").Value = 1.04
$ws.Range("N3").Value = 0.071
$ws.Range("O3").Value = 24.308

# Row 4 (large)
$ws.Range("C4").Value = 4.304
$ws.Range("D4").Value = 1.348
$ws.Range("E4").Value = 0.047
$ws.Range("F4").Value = 0.307
$ws.Range("G4").Value = 0.071
$ws.Range("H4").Value = 2.152
$ws.Range("I4").Value = 0.071
$ws.Range("J4").Value = 3.334
$ws.Range("K4").Value = 0.78
$ws.Range("L4").Value = 0.118
$ws.Range("M4").Value = 0.355
$ws.Range("N4").Value = 0.024
$ws.Range("O4").Value = 12.911

# Row 5 (huge)
$ws.Range("C5").Value = 0.024
$ws.Range("D5").Value = 0.213
$ws.Range("G5").Value = 0.024
$ws.Range("J5").Value = 0.024
$ws.Range("O5").Value = 0.285

# Row 6 (unknown_sz)
$ws.Range("C6").Value = 0.402
$ws.Range("D6").Value = 0.024
$ws.Range("E6").Value = 0.071
$ws.Range("F6").Value = 0.024
$ws.Range("G6").Value = 0.047
$ws.Range("H6").Value = 0.071
$ws.Range("I6").Value = 0.047
$ws.Range("J6").Value = 1.797
$ws.Range("K6").Value = 1.844
$ws.Range("L6").Value = 0.095
$ws.Range("M6").Value = 0.118
$ws.Range("N6").Value = 0.568
$ws.Range("O6").Value = 5.108

# Row 7 (COL_TOT)
$ws.Range("B7").Value = 0.071
$ws.Range("C7").Value = 21.803
$ws.Range("D7").Value = 1.94
$ws.Range("E7").Value = 0.236
$ws.Range("F7").Value = 1.253
$ws.Range("G7").Value = 0.497
$ws.Range("H7").Value = 4.375
$ws.Range("I7").Value = 0.639
$ws.Range("J7").Value = 41.003
$ws.Range("K7").Value = 17.757
$ws.Range("L7").Value = 5.226
$ws.Range("M7").Value = 2.601
$ws.Range("N7").Value = 2.602
$ws.Range("O7").Value = 100.003
